$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings: si 6 and si 9) ---
$ws.Range("A8").Characters(21, 2).Text = "50"
$ws.Range("C9").Characters(27, 9).Text = "12/8/2025"
$ws.Range("C9").Characters(47, 9).Text = "12/14/2025"

# --- Column H width (7.433768 -> 6.168446) ---
$ws.Columns("H").ColumnWidth = 6.168446

# --- Crime-stat table updates (rows 14-30) ---
# Row 14
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 54
$ws.Range("J15").Value = 44
$ws.Range("K15").Value = 22.727272727272
$ws.Range("L15").Value = 35
$ws.Range("M15").Value = 145.454545454545
$ws.Range("N15").Value = -31.645569620253
# Row 16
$ws.Range("C16").Value = 17
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 183.333333333333
$ws.Range("F16").Value = 43
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = 19.444444444444
$ws.Range("I16").Value = 601
$ws.Range("J16").Value = 604
$ws.Range("K16").Value = -0.496688741721
$ws.Range("L16").Value = -6.09375
$ws.Range("M16").Value = 39.767441860465
$ws.Range("N16").Value = -67.140513942044
# Row 17
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -31.25
$ws.Range("F17").Value = 64
$ws.Range("G17").Value = 80
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 993
$ws.Range("J17").Value = 992
$ws.Range("K17").Value = 0.100806451612
$ws.Range("L17").Value = 0.914634146341
$ws.Range("M17").Value = 139.277108433735
$ws.Range("N17").Value = -5.965909090909
# Row 18
$ws.Range("C18").Value = 16
$ws.Range("D18").Value = 8
$ws.Range("F18").Value = 37
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = 27.586206896551
$ws.Range("I18").Value = 429
$ws.Range("J18").Value = 345
$ws.Range("K18").Value = 24.347826086956
$ws.Range("L18").Value = 40.196078431372
$ws.Range("M18").Value = 122.279792746114
$ws.Range("N18").Value = -66.949152542372
# Row 19
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 111.111111111111
$ws.Range("F19").Value = 69
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = 25.454545454545
$ws.Range("I19").Value = 935
$ws.Range("J19").Value = 926
$ws.Range("K19").Value = 0.97192224622
$ws.Range("L19").Value = 30.952380952381
$ws.Range("M19").Value = 123.150357995227
$ws.Range("N19").Value = 27.557980900409
# Row 20
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 75
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 6.666666666666
$ws.Range("I20").Value = 255
$ws.Range("J20").Value = 227
$ws.Range("K20").Value = 12.334801762114
$ws.Range("L20").Value = -21.296296296296
$ws.Range("M20").Value = 117.948717948718
$ws.Range("N20").Value = -59.651898734177
# Row 21
$ws.Range("C21").Value = 71
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 61.363636363636
$ws.Range("F21").Value = 232
$ws.Range("G21").Value = 221
$ws.Range("H21").Value = 4.97737556561
$ws.Range("I21").Value = 3278
$ws.Range("J21").Value = 3152
$ws.Range("K21").Value = 3.997461928934
$ws.Range("L21").Value = 8.543046357615
$ws.Range("M21").Value = 103.60248447205
$ws.Range("N21").Value = -42.430628731998
# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("I22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -50
$ws.Range("K22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 16.666666666666
$ws.Range("I22").Value = 69
$ws.Range("J22").Value = 86
$ws.Range("K22").Value = -19.767441860465
$ws.Range("L22").Value = 6.153846153846
$ws.Range("M22").Value = 9.523809523809
# Row 23
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = 28.571428571428
$ws.Range("F23").Value = 40
$ws.Range("G23").Value = 35
$ws.Range("H23").Value = 14.285714285714
$ws.Range("I23").Value = 487
$ws.Range("J23").Value = 469
$ws.Range("K23").Value = 3.837953091684
$ws.Range("L23").Value = 3.837953091684
$ws.Range("M23").Value = 70.877192982456
# Row 24
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -17.647058823529
$ws.Range("F24").Value = 133
$ws.Range("G24").Value = 134
$ws.Range("H24").Value = -0.746268656716
$ws.Range("I24").Value = 1897
$ws.Range("J24").Value = 1660
$ws.Range("K24").Value = 14.277108433734
$ws.Range("L24").Value = 22.308188265635
$ws.Range("M24").Value = 38.974358974359
# Row 25
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 68
$ws.Range("H25").Value = -41.176470588235
$ws.Range("I25").Value = 764
$ws.Range("J25").Value = 737
$ws.Range("K25").Value = 3.663500678426
$ws.Range("L25").Value = 26.072607260726
# Row 26
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 96
$ws.Range("G26").Value = 78
$ws.Range("H26").Value = 23.076923076923
$ws.Range("I26").Value = 1228
$ws.Range("J26").Value = 1200
$ws.Range("K26").Value = 2.333333333333
$ws.Range("L26").Value = 13.284132841328
$ws.Range("M26").Value = 12.146118721461
# Row 27
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("I27").Value = 66
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = 15.78947368421
$ws.Range("L27").Value = 11.864406779661
# Row 28
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 13
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 62.5
$ws.Range("I28").Value = 116
$ws.Range("J28").Value = 154
$ws.Range("K28").Value = -24.675324675324
$ws.Range("L28").Value = 11.538461538461
# Row 29
$ws.Range("C29").Value = 2
$ws.Range("I29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F29").Value = 2
$ws.Range("I29").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -33.333333333333
$ws.Range("I29").Value = 35
$ws.Range("K29").Value = -31.372549019607
$ws.Range("L29").Value = 2.941176470588
$ws.Range("M29").Value = -40.677966101694
$ws.Range("N29").Value = -83.333333333333
# Row 30
$ws.Range("C30").Value = 2
$ws.Range("I30").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("F30").Value = 2
$ws.Range("I30").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -33.333333333333
$ws.Range("I30").Value = 29
$ws.Range("K30").Value = -36.95652173913
$ws.Range("L30").Value = -9.375
$ws.Range("M30").Value = -43.13725490196
$ws.Range("N30").Value = -84.656084656084
